$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text (masthead volume/number + reporting week date range) ---
$ws.Range("A8").Value = "Volume 30   Number  42"
$ws.Range("C9").Value = "Report Covering the Week  10/16/2023  Through  10/22/2023"

# --- Robbery row (16): was blank ("0"/"***.*" placeholders), now has real counts.
# Copy number formatting from a same-shaped donor cell (row 26, same column style),
# then overwrite with the real value, so the cell flips from text-placeholder to number. ---
$ws.Range("I26").Copy($ws.Range("C16"))
$ws.Range("C16").Value = 1
$ws.Range("J26").Copy($ws.Range("D16"))
$ws.Range("D16").Value = 1
$ws.Range("K26").Copy($ws.Range("E16"))
$ws.Range("E16").Value = 0

# --- Other Sex Crimes (27) / Shooting Vic. (28) / Shooting Inc. (29):
# D/E (and for row 27, C/D/E) flip from real numbers back to the "0"/"***.*" text placeholders.
# Copying directly from an existing placeholder cell (row 26) brings both the right format AND text. ---
$ws.Range("C26").Copy($ws.Range("C27"))
$ws.Range("C26").Copy($ws.Range("D27"))
$ws.Range("E26").Copy($ws.Range("E27"))
$ws.Range("C26").Copy($ws.Range("D28"))
$ws.Range("E26").Copy($ws.Range("E28"))
$ws.Range("C26").Copy($ws.Range("D29"))
$ws.Range("E26").Copy($ws.Range("E29"))

# --- Remaining cells: same type/format as before, only the numeric value changes. ---
$ws.Range("N15").Value = -27.272727272727
$ws.Range("I16").Value = 51
$ws.Range("J16").Value = 66
$ws.Range("K16").Value = -22.727272727272
$ws.Range("L16").Value = 168.421052631579
$ws.Range("M16").Value = -19.047619047619
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -75
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 6
$ws.Range("H17").Value = 16.666666666666
$ws.Range("I17").Value = 78
$ws.Range("J17").Value = 75
$ws.Range("K17").Value = 4
$ws.Range("L17").Value = 39.285714285714
$ws.Range("M17").Value = 73.333333333333
$ws.Range("N17").Value = -17.894736842105
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = 16.666666666666
$ws.Range("F18").Value = 24
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = 9.090909090909
$ws.Range("I18").Value = 247
$ws.Range("J18").Value = 216
$ws.Range("K18").Value = 14.351851851851
$ws.Range("L18").Value = 39.54802259887
$ws.Range("M18").Value = 22.277227722772
$ws.Range("N18").Value = -70.133010882708
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -30.76923076923
$ws.Range("F19").Value = 44
$ws.Range("G19").Value = 63
$ws.Range("H19").Value = -30.15873015873
$ws.Range("I19").Value = 510
$ws.Range("J19").Value = 526
$ws.Range("K19").Value = -3.041825095057
$ws.Range("L19").Value = 70.568561872909
$ws.Range("M19").Value = 73.469387755102
$ws.Range("N19").Value = 13.585746102449
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -25
$ws.Range("F20").Value = 28
$ws.Range("H20").Value = 154.545454545455
$ws.Range("I20").Value = 156
$ws.Range("J20").Value = 93
$ws.Range("K20").Value = 67.741935483871
$ws.Range("L20").Value = 173.684210526316
$ws.Range("M20").Value = 35.652173913043
$ws.Range("N20").Value = -94.099848714069
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 28
$ws.Range("E21").Value = -25
$ws.Range("F21").Value = 105
$ws.Range("G21").Value = 105
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 1052
$ws.Range("J21").Value = 978
$ws.Range("K21").Value = 7.566462167689
$ws.Range("L21").Value = 71.615008156606
$ws.Range("M21").Value = 45.303867403314
$ws.Range("N21").Value = -75.443510737628
$ws.Range("C24").Value = 8
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = -42.857142857142
$ws.Range("F24").Value = 49
$ws.Range("G24").Value = 65
$ws.Range("H24").Value = -24.615384615384
$ws.Range("I24").Value = 472
$ws.Range("J24").Value = 634
$ws.Range("K24").Value = -25.552050473186
$ws.Range("L24").Value = 7.762557077625
$ws.Range("M24").Value = 31.111111111111
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = -12.5
$ws.Range("F25").Value = 26
$ws.Range("G25").Value = 37
$ws.Range("H25").Value = -29.729729729729
$ws.Range("I25").Value = 188
$ws.Range("J25").Value = 202
$ws.Range("K25").Value = -6.930693069306
$ws.Range("L25").Value = 55.371900826446
$ws.Range("M25").Value = 35.251798561151
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 0
